$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 text values first (so "10" is the first brand-new shared string
# interned, matching the recorded edit chronology) ---
$ws.Cells.Item(11, 1).Value = "10"
$ws.Cells.Item(11, 2).Value = "10"

# --- Notes column (C2:C11): replace the repeated "notes" with distinct
# notes1..notes10 values, in row order ---
for ($i = 2; $i -le 11; $i++) {
    $n = $i - 1
    $ws.Cells.Item($i, 3).Value = "notes" + $n
}

# --- ID (A) / Label (B) columns for rows 2-10: increment 0..8 -> 1..9,
# keeping the values numeric even though the cell format is Text ("@").
# Toggling the format to General around the write avoids the "typed into a
# text cell" coercion, then restoring "@" keeps the original style index. ---
for ($i = 2; $i -le 10; $i++) {
    $n = $i - 1
    $ws.Cells.Item($i, 1).NumberFormat = "general"
    $ws.Cells.Item($i, 1).Value = $n
    $ws.Cells.Item($i, 1).NumberFormat = "@"

    $ws.Cells.Item($i, 2).NumberFormat = "general"
    $ws.Cells.Item($i, 2).Value = $n
    $ws.Cells.Item($i, 2).NumberFormat = "@"
}

# --- Update the active selection to match the saved view state ---
$ws.Range("E15").Select()
